$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sval data (regen to filter save games)
$data = @{
    2 = @{ B = 0.001754667048134761; C = 0.3375848360084654; D = 2938.103010863317;  E = 71517.89157740913; F = 0; G = 74456.3339277755 }
    3 = @{ B = 0.06328177979961902;  C = 9.226618575922256;  D = 157.8057217802531;   E = 71517.89157740913; F = 1; G = 71684.98719954511 }
    4 = @{ B = 0.1554434735375247;   C = 0.3375848360084654; D = 3.082599426703578;   E = 246.9852506941017; F = 1; G = 250.5608784303512 }
    5 = @{ B = 1.505614041169197;    C = 1.65323645889881;   D = 0.1529057820181812;  E = 0.4998867070740569; F = 1; G = 3.811642989160245 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
}
